# Fruta / hortaliza, semanal
# The data rows (2-15) get reshuffled: each target row receives the
# D (Fecha), I (Calidad), J (Volumen), K/L/M (Precios) and P (Precio $/Kg)
# values that used to live in a different source row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row number -> source row number (values to copy from)
$rowMap = @{
    2  = 15
    3  = 13
    4  = 8
    5  = 2
    6  = 7
    7  = 5
    8  = 11
    9  = 4
    10 = 12
    11 = 6
    12 = 10
    13 = 3
    14 = 9
    15 = 14
}

# Columns whose values move together with the row
$cols = @("D", "I", "J", "K", "L", "M", "P")

# 1) Snapshot the current ("before") values for every relevant cell
$snapshot = @{}
foreach ($row in 2..15) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

# 2) Write the values from the mapped source row into each target row
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $sourceData = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $sourceData[$col]
    }
}
